# feat: add 2022-Q1 data
#
# The existing "总计" (Total) sheet is renamed to "2022-Q1" and repopulated
# with that quarter's fund-holdings detail (same column layout used by the
# other quarterly sheets). A brand-new "总计" sheet is appended at the end
# (cloned from "2022-Q1" so it keeps the same sheet-level properties) with
# the refreshed roll-up table: old rows shifted down one, new 2022-Q1
# summary row inserted at the top.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the current last sheet ("总计") to "2022-Q1" and rewrite its
#    contents with the quarterly fund-holdings detail.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1.Name = "2022-Q1"

# Extend the existing header / index-column formatting (bold, centered,
# thin border) from the old 4-column layout onto the wider header row
# (B1:H1) and index column (A2:A13) so it keeps the look used throughout
# the workbook.
$q1.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$q1.Range("A2:A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Data rows: index, code, name, scale, total stock position, position share,
# held market value (亿元), position rank.
$q1Rows = @(
    @(0,  "012719", "华夏新兴经济一年持有期混合型证券投资基金A", "21.62", "91.19", "3.07", "0.6637", 9),
    @(1,  "159996", "国泰中证全指家用电器ETF",                     "17.91", "98.66", "1.89", "0.3385", 9),
    @(2,  "012421", "华夏优加生活混合A",                           "8.67",  "92.98", "2.76", "0.2393", 8),
    @(3,  "005063", "广发中证全指家用电器指数A",                   "9.91",  "94.24", "2.00", "0.1982", 10),
    @(4,  "003131", "国寿安保强国智造灵活配置混合",                "5.76",  "86.19", "1.92", "0.1106", 10),
    @(5,  "005064", "广发中证全指家用电器指数C",                   "5.40",  "94.24", "2.00", "0.1080", 10),
    @(6,  "004760", "国寿安保稳瑞混合A",                           "8.20",  "21.27", "0.71", "0.0582", 7),
    @(7,  "012720", "华夏新兴经济一年持有期混合型证券投资基金C", "1.04",  "91.19", "3.07", "0.0319", 9),
    @(8,  "010205", "国寿安保裕安混合A",                           "3.44",  "29.13", "0.75", "0.0258", 9),
    @(9,  "004761", "国寿安保稳瑞混合C",                           "2.52",  "21.27", "0.71", "0.0179", 7),
    @(10, "012422", "华夏优加生活混合C",                           "0.17",  "92.98", "2.76", "0.0047", 8),
    @(11, "010206", "国寿安保裕安混合C",                           "0.20",  "29.13", "0.75", "0.0015", 9)
)

$row = 2
foreach ($r in $q1Rows) {
    $q1.Cells.Item($row, 1).Value = $r[0]

    $q1.Cells.Item($row, 2).Value = "'" + $r[1]
    $q1.Cells.Item($row, 3).Value = $r[2]
    $q1.Cells.Item($row, 4).Value = "'" + $r[3]
    $q1.Cells.Item($row, 5).Value = "'" + $r[4]
    $q1.Cells.Item($row, 6).Value = "'" + $r[5]
    $q1.Cells.Item($row, 7).Value = "'" + $r[6]
    $q1.Cells.Item($row, 8).Value = $r[7]

    $row++
}

# ------------------------------------------------------------------
# 2. Clone "2022-Q1" (right after itself) to get a new sheet that keeps
#    the same sheet-level properties (outline/page setup), rename it to
#    "总计", wipe its cells, and write the refreshed roll-up table.
# ------------------------------------------------------------------
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item($wb.Worksheets.Count)
$total.Name = "总计"
$total.Cells.Clear()

# Re-use the same header / index-column formatting for the new sheet.
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 12, 1.8),
    @(1, "2021-Q4", 17, 2.15),
    @(2, "2021-Q3", 15, 2.88),
    @(3, "2021-Q2", 31, 5.4),
    @(4, "2021-Q1", 39, 4.5),
    @(5, "2020-Q4", 36, 7.14)
)

$row = 2
foreach ($r in $totalRows) {
    $total.Cells.Item($row, 1).Value = $r[0]

    $total.Cells.Item($row, 2).Value = $r[1]
    $total.Cells.Item($row, 3).Value = $r[2]
    $total.Cells.Item($row, 4).Value = $r[3]

    $row++
}

$wb.Save()
